$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update counts in column C
$ws.Range("C2").Value = 19
$ws.Range("C3").Value = 20
$ws.Range("C6").Value = 17
$ws.Range("C7").Value = 10
$ws.Range("C8").Value = 14
$ws.Range("C9").Value = 16
$ws.Range("C10").Value = 22
$ws.Range("C11").Value = 25
$ws.Range("C12").Value = 18
$ws.Range("C14").Value = 13
$ws.Range("C15").Value = 23
$ws.Range("C16").Value = 16
$ws.Range("C17").Value = 21
$ws.Range("C18").Value = 15

# Update text label in column B
$ws.Range("B17").Value = "<hin>"
